$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item(1).Name = "services"
$wb.Worksheets.Item(2).Name = "other ports"

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Move the cursor on sheet1 (services) to A46
$ws1.Range("A46").Select()

# Populate sheet2 ("other ports") with the new table
$ws2.Range("A1").Value = "caGrid 1.0 Training"
$ws2.Range("A2").Value = "Host:Port"
$ws2.Range("B2").Value = "user"
$ws2.Range("C2").Value = "Protocol"
$ws2.Range("D2").Value = "Description"
$ws2.Range("A3").Value = "usage.cagrid.org:55555"
$ws2.Range("B3").Value = "introduce"
$ws2.Range("C3").Value = "UDP"
$ws2.Range("D3").Value = "introduce stats collecting"

$ws2.Range("D3").Select()
